$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking price text (e.g. "605.93") that the
# COM layer would otherwise auto-convert to a real number. Force them to stay
# plain text (matching the original inlineStr cells, with no cell style
# applied) by switching to a text number format before assignment, then
# resetting the style back to Normal so no stray style index is left behind.
$priceUpdates = @{
    "D2"  = "64.341.26"
    "D3"  = "3.145.06"
    "D5"  = "605.93"
    "D6"  = "149.47"
    "D8"  = "3.141.79"
    "D11" = "5.60"
    "D14" = "36.87"
    "D15" = "3.659.70"
    "D16" = "64.381.54"
    "D18" = "3.143.79"
    "D20" = "482.51"
    "D21" = "14.64"
    "D22" = "0.713"
    "D23" = "7.77"
    "D24" = "13.83"
    "D25" = "84.19"
    "D28" = "8.56"
    "D30" = "0.125"
    "D31" = "7.00"
    "D32" = "2.72"
    "D33" = "1.00"
    "D34" = "26.75"
    "D40" = "452.66"
    "D43" = "8.49"
    "D44" = "2.888.36"
    "D45" = "0.273"
    "D46" = "2.32"
    "D47" = "26.76"
    "D48" = "0.999"
    "D51" = "119.59"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}

# Remaining text cells (percentages in column E, plus the renamed coin row)
# are not number-like, so they can be assigned directly without any
# auto-conversion risk.
$textUpdates = @{
    "E2"  = "  -2.91%  "
    "E3"  = "  -1.81%  "
    "E4"  = "  +0.09%  "
    "E5"  = "  -0.23%  "
    "E6"  = "  -4.19%  "
    "E7"  = "  +0.06%  "
    "E8"  = "  -1.86%  "
    "E9"  = "  -3.22%  "
    "E10" = "  -4.60%  "
    "E11" = "  -0.97%  "
    "E12" = "  -4.60%  "
    "E13" = "  -2.93%  "
    "E14" = "  -3.83%  "
    "E15" = "  -1.85%  "
    "E16" = "  -3.06%  "
    "E17" = "  +0.10%  "
    "E18" = "  -1.91%  "
    "E19" = "  -3.98%  "
    "E20" = "  -4.67%  "
    "E21" = "  -4.32%  "
    "E22" = "  -2.19%  "
    "E23" = "  -2.72%  "
    "E24" = "  -5.17%  "
    "E25" = "  -1.06%  "
    "E26" = "  +0.12%  "
    "E27" = "  -1.76%  "
    "E28" = "  -5.16%  "
    "E29" = "  -4.11%  "
    "E30" = "  -2.54%  "
    "E31" = "  +0.76%  "
    "E32" = "  -6.76%  "
    "E33" = "  -0.17%  "
    "E34" = "  -5.17%  "
    "E35" = "  -5.15%  "
    "E36" = "  -4.94%  "
    "E37" = "  -1.59%  "
    "E38" = "  +7.22%  "
    "E39" = "  -1.26%  "
    "E40" = "  -9.56%  "
    "E41" = "  -4.02%  "
    "E42" = "  -5.43%  "
    "E43" = "  -2.46%  "
    "E44" = "  -0.79%  "
    "E45" = "  -7.57%  "
    "E46" = "  -4.52%  "
    "E48" = "  -0.02%  "
    "E49" = "  -0.90%  "
    "E50" = "  -3.02%  "
    "B51" = "Monero"
    "C51" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "E51" = "  -2.04%  "
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}
